$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 961.1667
$ws.Cells.Item(33, 9).Value = 653.6
$ws.Cells.Item(33, 11).Value = 653.6
$ws.Cells.Item(33, 13).Value = -424.6
$ws.Cells.Item(57, 8).Value = 102854.5
$ws.Cells.Item(57, 9).Value = 45709
$ws.Cells.Item(57, 10).Value = 160000
$ws.Cells.Item(57, 11).Value = 137127
$ws.Cells.Item(57, 12).Value = 480000
$ws.Cells.Item(57, 13).Value = -136628
$ws.Cells.Item(57, 14).Value = -480998
$ws.Cells.Item(61, 8).Value = 0
$ws.Cells.Item(61, 9).Value = 0
$ws.Cells.Item(61, 11).Value = 0
$ws.Cells.Item(61, 13).ClearContents()
$ws.Cells.Item(62, 8).Value = 9317.4
$ws.Cells.Item(62, 9).Value = 8081.6665
$ws.Cells.Item(62, 10).Value = 11171
$ws.Cells.Item(62, 11).Value = 8081.6665
$ws.Cells.Item(62, 12).Value = 11171
$ws.Cells.Item(62, 13).Value = -7457.6665
$ws.Cells.Item(62, 14).Value = -12419
$ws.Cells.Item(65, 8).Value = 9317.4
$ws.Cells.Item(65, 9).Value = 8081.6665
$ws.Cells.Item(65, 10).Value = 11171
$ws.Cells.Item(65, 11).Value = 40408.3325
$ws.Cells.Item(65, 12).Value = 55855
$ws.Cells.Item(65, 13).Value = -37288.3325
$ws.Cells.Item(65, 14).Value = -62095
$ws.Cells.Item(86, 8).Value = 132426650
$ws.Cells.Item(86, 10).Value = 250206990
$ws.Cells.Item(86, 12).Value = 250206990
$ws.Cells.Item(86, 14).Value = -250209236
$ws.Cells.Item(89, 8).Value = 132426650
$ws.Cells.Item(89, 9).Value = 83351496
$ws.Cells.Item(89, 10).Value = 250206990
$ws.Cells.Item(89, 12).Value = 1251034950
$ws.Cells.Item(89, 14).Value = -1251046182
$ws.Cells.Item(98, 8).Value = 2956
$ws.Cells.Item(98, 9).Value = 2372.1428
$ws.Cells.Item(98, 11).Value = 2372.1428
$ws.Cells.Item(98, 13).Value = -874.1428000000001
$ws.Cells.Item(99, 8).Value = 885.25
$ws.Cells.Item(99, 9).Value = 1120
$ws.Cells.Item(99, 10).Value = 181
$ws.Cells.Item(99, 11).Value = 3360
$ws.Cells.Item(99, 12).Value = 543
$ws.Cells.Item(99, 13).Value = -1862
$ws.Cells.Item(99, 14).Value = -3539
$ws.Cells.Item(111, 8).Value = 3198.6
$ws.Cells.Item(111, 9).Value = 3500
$ws.Cells.Item(111, 11).Value = 10500
$ws.Cells.Item(111, 13).Value = -7433
$ws.Cells.Item(122, 8).Value = 2956
$ws.Cells.Item(122, 9).Value = 2372.1428
$ws.Cells.Item(122, 11).Value = 7116.428400000001
$ws.Cells.Item(122, 13).Value = -4666.428400000001
$ws.Cells.Item(132, 8).Value = 1365.6428
$ws.Cells.Item(132, 9).Value = 1382.6538
$ws.Cells.Item(132, 11).Value = 4147.9614
$ws.Cells.Item(132, 13).Value = -1617.9614
$ws.Cells.Item(138, 8).Value = 3408.4
$ws.Cells.Item(138, 9).Value = 1349.76
$ws.Cells.Item(138, 10).Value = 4094.6133
$ws.Cells.Item(138, 11).Value = 4049.28
$ws.Cells.Item(138, 12).Value = 12283.8399
$ws.Cells.Item(138, 13).Value = 1090.72
$ws.Cells.Item(138, 14).Value = -22563.8399
$ws.Cells.Item(141, 8).Value = 444.33334
$ws.Cells.Item(141, 9).Value = 444.33334
$ws.Cells.Item(141, 11).Value = 1333.00002
$ws.Cells.Item(141, 13).Value = 3846.99998

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 20590636
$ws.Cells.Item(32, 9).Value = 22380452
$ws.Cells.Item(32, 10).Value = 7749.75
$ws.Cells.Item(32, 11).Value = 22380452
$ws.Cells.Item(32, 12).Value = 7749.75
$ws.Cells.Item(32, 13).Value = -22380165
$ws.Cells.Item(32, 14).Value = -8323.75
$ws.Cells.Item(45, 8).Value = 5265.3335
$ws.Cells.Item(45, 9).Value = 5107.4443
$ws.Cells.Item(45, 10).Value = 5502.1665
$ws.Cells.Item(45, 11).Value = 5107.4443
$ws.Cells.Item(45, 12).Value = 5502.1665
$ws.Cells.Item(45, 13).Value = -4730.4443
$ws.Cells.Item(45, 14).Value = -6256.1665
$ws.Cells.Item(102, 8).Value = 2346.6155
$ws.Cells.Item(102, 9).Value = 1001.5
$ws.Cells.Item(102, 11).Value = 1001.5
$ws.Cells.Item(102, 13).Value = 620.5
$ws.Cells.Item(113, 8).Value = 49999
$ws.Cells.Item(113, 10).Value = 49999
$ws.Cells.Item(113, 12).Value = 49999
$ws.Cells.Item(113, 14).Value = -58677
$ws.Cells.Item(132, 8).Value = 2494.9185
$ws.Cells.Item(132, 9).Value = 2225.3684
$ws.Cells.Item(132, 10).Value = 3426.0908
$ws.Cells.Item(132, 11).Value = 6676.1052
$ws.Cells.Item(132, 12).Value = 10278.2724
$ws.Cells.Item(132, 13).Value = -4146.1052
$ws.Cells.Item(132, 14).Value = -15338.2724
$ws.Cells.Item(134, 8).Value = 159998.5
$ws.Cells.Item(134, 10).Value = 159998.5
$ws.Cells.Item(134, 12).Value = 159998.5
$ws.Cells.Item(134, 14).Value = -170138.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 1658.5
$ws.Cells.Item(20, 9).Value = 1032.1875
$ws.Cells.Item(20, 11).Value = 1032.1875
$ws.Cells.Item(20, 13).Value = -785.1875
$ws.Cells.Item(86, 8).Value = 5000
$ws.Cells.Item(86, 9).Value = 5000
$ws.Cells.Item(86, 11).Value = 5000
$ws.Cells.Item(86, 13).Value = -3877
$ws.Cells.Item(89, 8).Value = 5000
$ws.Cells.Item(89, 9).Value = 5000
$ws.Cells.Item(89, 11).Value = 25000
$ws.Cells.Item(89, 13).Value = -19384
$ws.Cells.Item(94, 8).Value = 1077.9642
$ws.Cells.Item(94, 9).Value = 832
$ws.Cells.Item(94, 10).Value = 1458.091
$ws.Cells.Item(94, 11).Value = 832
$ws.Cells.Item(94, 12).Value = 1458.091
$ws.Cells.Item(94, 13).Value = -381
$ws.Cells.Item(94, 14).Value = -2360.091
$ws.Cells.Item(105, 8).Value = 1971.4
$ws.Cells.Item(105, 9).Value = 1736.3077
$ws.Cells.Item(105, 11).Value = 1736.3077
$ws.Cells.Item(105, 13).Value = 10.69229999999993
$ws.Cells.Item(107, 8).Value = 2077.4546
$ws.Cells.Item(107, 9).Value = 1416.2858
$ws.Cells.Item(107, 11).Value = 1416.2858
$ws.Cells.Item(107, 13).Value = 503.7141999999999
$ws.Cells.Item(125, 8).Value = 122000
$ws.Cells.Item(125, 10).Value = 122000
$ws.Cells.Item(125, 12).Value = 122000
$ws.Cells.Item(125, 14).Value = -131840
$ws.Cells.Item(133, 8).Value = 115250.75
$ws.Cells.Item(133, 10).Value = 115250.75
$ws.Cells.Item(133, 12).Value = 115250.75
$ws.Cells.Item(133, 14).Value = -125370.75
$ws.Cells.Item(134, 8).Value = 4205269.5
$ws.Cells.Item(134, 9).Value = 5496535.5
$ws.Cells.Item(134, 11).Value = 16489606.5
$ws.Cells.Item(134, 13).Value = -16487071.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(109, 8).Value = 26415
$ws.Cells.Item(109, 10).Value = 26415
$ws.Cells.Item(109, 12).Value = 26415
$ws.Cells.Item(109, 14).Value = -28495
$ws.Cells.Item(112, 8).Value = 80301
$ws.Cells.Item(112, 10).Value = 80301
$ws.Cells.Item(112, 12).Value = 80301
$ws.Cells.Item(112, 14).Value = -83255

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(4, 8).Value = 114067040
$ws.Cells.Item(4, 9).Value = 161425100
$ws.Cells.Item(4, 11).Value = 484275300
$ws.Cells.Item(4, 13).Value = -484275188
$ws.Cells.Item(32, 8).Value = 125000440
$ws.Cells.Item(32, 9).Value = 499.85715
$ws.Cells.Item(32, 10).Value = 1000000000
$ws.Cells.Item(32, 11).Value = 1499.57145
$ws.Cells.Item(32, 12).Value = 3000000000
$ws.Cells.Item(32, 13).Value = -1216.57145
$ws.Cells.Item(32, 14).Value = -3000000566
$ws.Cells.Item(54, 8).Value = 0
$ws.Cells.Item(54, 10).Value = 0
$ws.Cells.Item(54, 12).Value = 0
$ws.Cells.Item(54, 14).ClearContents()
$ws.Cells.Item(120, 8).Value = 21825.309
$ws.Cells.Item(120, 9).Value = 17746
$ws.Cells.Item(120, 10).Value = 24374.875
$ws.Cells.Item(120, 11).Value = 53238
$ws.Cells.Item(120, 12).Value = 73124.625
$ws.Cells.Item(120, 13).Value = -48400
$ws.Cells.Item(120, 14).Value = -82800.625
$ws.Cells.Item(122, 8).Value = 537045.4399999999
$ws.Cells.Item(122, 10).Value = 1162549
$ws.Cells.Item(122, 12).Value = 10462941
$ws.Cells.Item(122, 14).Value = -10467841
$ws.Cells.Item(128, 8).Value = 165000
$ws.Cells.Item(128, 9).Value = 165000
$ws.Cells.Item(128, 11).Value = 495000
$ws.Cells.Item(128, 13).Value = -490020

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 3483.2856
$ws.Cells.Item(102, 9).Value = 3397.25
$ws.Cells.Item(102, 11).Value = 3397.25
$ws.Cells.Item(102, 13).Value = -1775.25
$ws.Cells.Item(122, 8).Value = 3115.0667
$ws.Cells.Item(122, 9).Value = 3803
$ws.Cells.Item(122, 10).Value = 2328.8572
$ws.Cells.Item(122, 11).Value = 11409
$ws.Cells.Item(122, 12).Value = 6986.571599999999
$ws.Cells.Item(122, 13).Value = -8959
$ws.Cells.Item(122, 14).Value = -11886.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(46, 8).Value = 4205.136
$ws.Cells.Item(46, 9).Value = 996
$ws.Cells.Item(46, 10).Value = 4918.278
$ws.Cells.Item(46, 11).Value = 996
$ws.Cells.Item(46, 12).Value = 4918.278
$ws.Cells.Item(46, 13).Value = -808
$ws.Cells.Item(46, 14).Value = -5294.278
$ws.Cells.Item(68, 8).Value = 5209.68
$ws.Cells.Item(68, 9).Value = 4655.8823
$ws.Cells.Item(68, 10).Value = 6386.5
$ws.Cells.Item(68, 11).Value = 4655.8823
$ws.Cells.Item(68, 12).Value = 6386.5
$ws.Cells.Item(68, 13).Value = -3906.8823
$ws.Cells.Item(68, 14).Value = -7884.5
$ws.Cells.Item(71, 8).Value = 5209.68
$ws.Cells.Item(71, 9).Value = 4655.8823
$ws.Cells.Item(71, 10).Value = 6386.5
$ws.Cells.Item(71, 11).Value = 23279.4115
$ws.Cells.Item(71, 12).Value = 31932.5
$ws.Cells.Item(71, 13).Value = -19535.4115
$ws.Cells.Item(71, 14).Value = -39420.5
$ws.Cells.Item(122, 8).Value = 10999.2
$ws.Cells.Item(122, 9).Value = 6249.25
$ws.Cells.Item(122, 11).Value = 18747.75
$ws.Cells.Item(122, 13).Value = -16297.75
$ws.Cells.Item(132, 8).Value = 6759.522
$ws.Cells.Item(132, 9).Value = 6427.095
$ws.Cells.Item(132, 10).Value = 10250
$ws.Cells.Item(132, 11).Value = 19281.285
$ws.Cells.Item(132, 12).Value = 30750
$ws.Cells.Item(132, 13).Value = -16751.285
$ws.Cells.Item(132, 14).Value = -35810

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4820.1665
$ws.Cells.Item(62, 9).Value = 3240.3333
$ws.Cells.Item(62, 10).Value = 6400
$ws.Cells.Item(62, 11).Value = 3240.3333
$ws.Cells.Item(62, 12).Value = 6400
$ws.Cells.Item(62, 13).Value = -2616.3333
$ws.Cells.Item(62, 14).Value = -7648
$ws.Cells.Item(65, 8).Value = 4820.1665
$ws.Cells.Item(65, 9).Value = 3240.3333
$ws.Cells.Item(65, 10).Value = 6400
$ws.Cells.Item(65, 11).Value = 16201.6665
$ws.Cells.Item(65, 12).Value = 32000
$ws.Cells.Item(65, 13).Value = -13081.6665
$ws.Cells.Item(65, 14).Value = -38240
$ws.Cells.Item(81, 8).Value = 3672.2
$ws.Cells.Item(81, 9).Value = 2957.8462
$ws.Cells.Item(81, 10).Value = 4998.857
$ws.Cells.Item(81, 11).Value = 5915.6924
$ws.Cells.Item(81, 12).Value = 9997.714
$ws.Cells.Item(81, 13).Value = -4854.6924
$ws.Cells.Item(81, 14).Value = -12119.714
$ws.Cells.Item(84, 8).Value = 3672.2
$ws.Cells.Item(84, 9).Value = 2957.8462
$ws.Cells.Item(84, 10).Value = 4998.857
$ws.Cells.Item(84, 11).Value = 29578.462
$ws.Cells.Item(84, 12).Value = 49988.57
$ws.Cells.Item(84, 13).Value = -24274.462
$ws.Cells.Item(84, 14).Value = -60596.57
$ws.Cells.Item(125, 8).Value = 111995
$ws.Cells.Item(125, 10).Value = 111995
$ws.Cells.Item(125, 12).Value = 111995
$ws.Cells.Item(125, 14).Value = -121835
$ws.Cells.Item(132, 8).Value = 1860.7847
$ws.Cells.Item(132, 9).Value = 1758.2623
$ws.Cells.Item(132, 11).Value = 5274.7869
$ws.Cells.Item(132, 13).Value = -2744.7869
